$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "42.825.50"
$ws.Range("E2").Value = "  -1.61%  "

# Row 3
$ws.Range("D3").Value = "2.335.85"
$ws.Range("E3").Value = "  +0.68%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.03%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "306.52"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.75%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "100.36"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.25%  "

# Row 7
$ws.Range("E7").Value = "  -5.20%  "

# Row 8
$ws.Range("E8").Value = "  -0.06%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.510"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.63%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.99"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.13%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "52.16"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.30%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0800"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.23%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.81"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.44%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.79"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +5.06%  "

# Row 16
$ws.Range("D16").Value = "2.275.84"
$ws.Range("E16").Value = "  -1.57%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.796"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.19%  "

# Row 18
$ws.Range("D18").Value = "42.732.12"
$ws.Range("E18").Value = "  -1.59%  "

# Row 19
$ws.Range("E19").Value = "  +1.61%  "

# Row 20
$ws.Range("D20").Value = "0.0₃0906"
$ws.Range("E20").Value = "  -2.42%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.64"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -7.01%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.74"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.01%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "236.68"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.38%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.00"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.50%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.57"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.44%  "

# Row 26
$ws.Range("E26").Value = "  +0.41%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.97"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.76%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.27"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +7.38%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "34.97"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -7.23%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.37"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.10%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "159.47"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.62%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.00"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.06%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.12"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.15%  "

# Row 34
$ws.Range("B34").Value = "RenderToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.61"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +6.44%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "17.42"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.16%  "

# Row 36
$ws.Range("B36").Value = "WEMIXToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.45"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.03%  "

# Row 37
$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0728"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.72%  "

# Row 38
$ws.Range("E38").Value = "  -4.76%  "

# Row 39
$ws.Range("E39").Value = "  -0.59%  "

# Row 40
$ws.Range("E40").Value = "  -3.75%  "

# Row 41
$ws.Range("E41").Value = "  -3.47%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.36"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.68%  "

# Row 43
$ws.Range("D43").Value = "2.016.19"
$ws.Range("E43").Value = "  +1.71%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0285"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.94%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "18.71"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.54%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.31"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.32%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.93"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.82%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "56.10"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.13%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.91"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.28%  "

# Row 50
$ws.Range("D50").Value = "2.560.38"
$ws.Range("E50").Value = "  +0.56%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.64"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.35%  "
